$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the weekly Fruta/Hortaliza (Damasco) update: 3 new rows inserted at the top of the block
# (rows 179-181, new date 2023-12-07) and the remaining rows shifted/re-sorted down, with 3
# additional historical rows appended at the end (rows 200-202). Overwrite the whole A179:T202
# block explicitly to land on the exact post-edit layout, then extend the sheet dimension.

# Row 179
$ws.Cells.Item(179, 1).Value = 9
$ws.Cells.Item(179, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(179, 3).Value = "Metropolitana"
$ws.Cells.Item(179, 4).Value = 45267
$ws.Cells.Item(179, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(179, 5).Value = 13
$ws.Cells.Item(179, 6).Value = "Fruta"
$ws.Cells.Item(179, 7).Value = 100103
$ws.Cells.Item(179, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(179, 9).Value = 100103003
$ws.Cells.Item(179, 10).Value = "Damasco"
$ws.Cells.Item(179, 11).Value = "Castle Brite"
$ws.Cells.Item(179, 12).Value = "Especial"
$ws.Cells.Item(179, 13).Value = 80
$ws.Cells.Item(179, 14).Value = 20000
$ws.Cells.Item(179, 15).Value = 20000
$ws.Cells.Item(179, 16).Value = 20000
$ws.Cells.Item(179, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(179, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(179, 19).Value = 2000
$ws.Cells.Item(179, 20).Value = 10

# Row 180
$ws.Cells.Item(180, 1).Value = 9
$ws.Cells.Item(180, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(180, 3).Value = "Metropolitana"
$ws.Cells.Item(180, 4).Value = 45267
$ws.Cells.Item(180, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(180, 5).Value = 13
$ws.Cells.Item(180, 6).Value = "Fruta"
$ws.Cells.Item(180, 7).Value = 100103
$ws.Cells.Item(180, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(180, 9).Value = 100103003
$ws.Cells.Item(180, 10).Value = "Damasco"
$ws.Cells.Item(180, 11).Value = "Castle Brite"
$ws.Cells.Item(180, 12).Value = "Primera"
$ws.Cells.Item(180, 13).Value = 100
$ws.Cells.Item(180, 14).Value = 16000
$ws.Cells.Item(180, 15).Value = 16000
$ws.Cells.Item(180, 16).Value = 16000
$ws.Cells.Item(180, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(180, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(180, 19).Value = 1600
$ws.Cells.Item(180, 20).Value = 10

# Row 181
$ws.Cells.Item(181, 1).Value = 9
$ws.Cells.Item(181, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(181, 3).Value = "Metropolitana"
$ws.Cells.Item(181, 4).Value = 45267
$ws.Cells.Item(181, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(181, 5).Value = 13
$ws.Cells.Item(181, 6).Value = "Fruta"
$ws.Cells.Item(181, 7).Value = 100103
$ws.Cells.Item(181, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(181, 9).Value = 100103003
$ws.Cells.Item(181, 10).Value = "Damasco"
$ws.Cells.Item(181, 11).Value = "Castle Brite"
$ws.Cells.Item(181, 12).Value = "Segunda"
$ws.Cells.Item(181, 13).Value = 70
$ws.Cells.Item(181, 14).Value = 12000
$ws.Cells.Item(181, 15).Value = 12000
$ws.Cells.Item(181, 16).Value = 12000
$ws.Cells.Item(181, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(181, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(181, 19).Value = 1200
$ws.Cells.Item(181, 20).Value = 10

# Row 182
$ws.Cells.Item(182, 1).Value = 9
$ws.Cells.Item(182, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(182, 3).Value = "Metropolitana"
$ws.Cells.Item(182, 4).Value = 44187
$ws.Cells.Item(182, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(182, 5).Value = 13
$ws.Cells.Item(182, 6).Value = "Fruta"
$ws.Cells.Item(182, 7).Value = 100103
$ws.Cells.Item(182, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(182, 9).Value = 100103003
$ws.Cells.Item(182, 10).Value = "Damasco"
$ws.Cells.Item(182, 11).Value = "Patterson"
$ws.Cells.Item(182, 12).Value = "Primera"
$ws.Cells.Item(182, 13).Value = 80
$ws.Cells.Item(182, 14).Value = 15000
$ws.Cells.Item(182, 15).Value = 15000
$ws.Cells.Item(182, 16).Value = 15000
$ws.Cells.Item(182, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(182, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(182, 19).Value = 1000
$ws.Cells.Item(182, 20).Value = 15

# Row 183
$ws.Cells.Item(183, 1).Value = 9
$ws.Cells.Item(183, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(183, 3).Value = "Metropolitana"
$ws.Cells.Item(183, 4).Value = 44187
$ws.Cells.Item(183, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(183, 5).Value = 13
$ws.Cells.Item(183, 6).Value = "Fruta"
$ws.Cells.Item(183, 7).Value = 100103
$ws.Cells.Item(183, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(183, 9).Value = 100103003
$ws.Cells.Item(183, 10).Value = "Damasco"
$ws.Cells.Item(183, 11).Value = "Patterson"
$ws.Cells.Item(183, 12).Value = "Segunda"
$ws.Cells.Item(183, 13).Value = 95
$ws.Cells.Item(183, 14).Value = 13500
$ws.Cells.Item(183, 15).Value = 13500
$ws.Cells.Item(183, 16).Value = 13500
$ws.Cells.Item(183, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(183, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(183, 19).Value = 900
$ws.Cells.Item(183, 20).Value = 15

# Row 184
$ws.Cells.Item(184, 1).Value = 9
$ws.Cells.Item(184, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(184, 3).Value = "Metropolitana"
$ws.Cells.Item(184, 4).Value = 44187
$ws.Cells.Item(184, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(184, 5).Value = 13
$ws.Cells.Item(184, 6).Value = "Fruta"
$ws.Cells.Item(184, 7).Value = 100103
$ws.Cells.Item(184, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(184, 9).Value = 100103003
$ws.Cells.Item(184, 10).Value = "Damasco"
$ws.Cells.Item(184, 11).Value = "Patterson"
$ws.Cells.Item(184, 12).Value = "Tercera"
$ws.Cells.Item(184, 13).Value = 120
$ws.Cells.Item(184, 14).Value = 12000
$ws.Cells.Item(184, 15).Value = 12000
$ws.Cells.Item(184, 16).Value = 12000
$ws.Cells.Item(184, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(184, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(184, 19).Value = 800
$ws.Cells.Item(184, 20).Value = 15

# Row 185
$ws.Cells.Item(185, 1).Value = 9
$ws.Cells.Item(185, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(185, 3).Value = "Metropolitana"
$ws.Cells.Item(185, 4).Value = 44525
$ws.Cells.Item(185, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(185, 5).Value = 13
$ws.Cells.Item(185, 6).Value = "Fruta"
$ws.Cells.Item(185, 7).Value = 100103
$ws.Cells.Item(185, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(185, 9).Value = 100103003
$ws.Cells.Item(185, 10).Value = "Damasco"
$ws.Cells.Item(185, 11).Value = "Castle Brite"
$ws.Cells.Item(185, 12).Value = "Especial"
$ws.Cells.Item(185, 13).Value = 300
$ws.Cells.Item(185, 14).Value = 25200
$ws.Cells.Item(185, 15).Value = 25200
$ws.Cells.Item(185, 16).Value = 25200
$ws.Cells.Item(185, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(185, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(185, 19).Value = 1400
$ws.Cells.Item(185, 20).Value = 18

# Row 186
$ws.Cells.Item(186, 1).Value = 9
$ws.Cells.Item(186, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(186, 3).Value = "Metropolitana"
$ws.Cells.Item(186, 4).Value = 44525
$ws.Cells.Item(186, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(186, 5).Value = 13
$ws.Cells.Item(186, 6).Value = "Fruta"
$ws.Cells.Item(186, 7).Value = 100103
$ws.Cells.Item(186, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(186, 9).Value = 100103003
$ws.Cells.Item(186, 10).Value = "Damasco"
$ws.Cells.Item(186, 11).Value = "Castle Brite"
$ws.Cells.Item(186, 12).Value = "Primera"
$ws.Cells.Item(186, 13).Value = 250
$ws.Cells.Item(186, 14).Value = 21600
$ws.Cells.Item(186, 15).Value = 21600
$ws.Cells.Item(186, 16).Value = 21600
$ws.Cells.Item(186, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(186, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(186, 19).Value = 1200
$ws.Cells.Item(186, 20).Value = 18

# Row 187
$ws.Cells.Item(187, 1).Value = 9
$ws.Cells.Item(187, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(187, 3).Value = "Metropolitana"
$ws.Cells.Item(187, 4).Value = 44559
$ws.Cells.Item(187, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(187, 5).Value = 13
$ws.Cells.Item(187, 6).Value = "Fruta"
$ws.Cells.Item(187, 7).Value = 100103
$ws.Cells.Item(187, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(187, 9).Value = 100103003
$ws.Cells.Item(187, 10).Value = "Damasco"
$ws.Cells.Item(187, 11).Value = "Helena"
$ws.Cells.Item(187, 12).Value = "Especial"
$ws.Cells.Item(187, 13).Value = 310
$ws.Cells.Item(187, 14).Value = 15000
$ws.Cells.Item(187, 15).Value = 15000
$ws.Cells.Item(187, 16).Value = 15000
$ws.Cells.Item(187, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(187, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(187, 19).Value = 1000
$ws.Cells.Item(187, 20).Value = 15

# Row 188
$ws.Cells.Item(188, 1).Value = 9
$ws.Cells.Item(188, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(188, 3).Value = "Metropolitana"
$ws.Cells.Item(188, 4).Value = 44559
$ws.Cells.Item(188, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(188, 5).Value = 13
$ws.Cells.Item(188, 6).Value = "Fruta"
$ws.Cells.Item(188, 7).Value = 100103
$ws.Cells.Item(188, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(188, 9).Value = 100103003
$ws.Cells.Item(188, 10).Value = "Damasco"
$ws.Cells.Item(188, 11).Value = "Helena"
$ws.Cells.Item(188, 12).Value = "Primera"
$ws.Cells.Item(188, 13).Value = 350
$ws.Cells.Item(188, 14).Value = 12000
$ws.Cells.Item(188, 15).Value = 12000
$ws.Cells.Item(188, 16).Value = 12000
$ws.Cells.Item(188, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(188, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(188, 19).Value = 800
$ws.Cells.Item(188, 20).Value = 15

# Row 189
$ws.Cells.Item(189, 1).Value = 9
$ws.Cells.Item(189, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(189, 3).Value = "Metropolitana"
$ws.Cells.Item(189, 4).Value = 44559
$ws.Cells.Item(189, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(189, 5).Value = 13
$ws.Cells.Item(189, 6).Value = "Fruta"
$ws.Cells.Item(189, 7).Value = 100103
$ws.Cells.Item(189, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(189, 9).Value = 100103003
$ws.Cells.Item(189, 10).Value = "Damasco"
$ws.Cells.Item(189, 11).Value = "Patterson"
$ws.Cells.Item(189, 12).Value = "Especial"
$ws.Cells.Item(189, 13).Value = 300
$ws.Cells.Item(189, 14).Value = 18000
$ws.Cells.Item(189, 15).Value = 18000
$ws.Cells.Item(189, 16).Value = 18000
$ws.Cells.Item(189, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(189, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(189, 19).Value = 1000
$ws.Cells.Item(189, 20).Value = 18

# Row 190
$ws.Cells.Item(190, 1).Value = 9
$ws.Cells.Item(190, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(190, 3).Value = "Metropolitana"
$ws.Cells.Item(190, 4).Value = 44559
$ws.Cells.Item(190, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(190, 5).Value = 13
$ws.Cells.Item(190, 6).Value = "Fruta"
$ws.Cells.Item(190, 7).Value = 100103
$ws.Cells.Item(190, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(190, 9).Value = 100103003
$ws.Cells.Item(190, 10).Value = "Damasco"
$ws.Cells.Item(190, 11).Value = "Patterson"
$ws.Cells.Item(190, 12).Value = "Primera"
$ws.Cells.Item(190, 13).Value = 380
$ws.Cells.Item(190, 14).Value = 14400
$ws.Cells.Item(190, 15).Value = 14400
$ws.Cells.Item(190, 16).Value = 14400
$ws.Cells.Item(190, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(190, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(190, 19).Value = 800
$ws.Cells.Item(190, 20).Value = 18

# Row 191
$ws.Cells.Item(191, 1).Value = 9
$ws.Cells.Item(191, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(191, 3).Value = "Metropolitana"
$ws.Cells.Item(191, 4).Value = 44558
$ws.Cells.Item(191, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(191, 5).Value = 13
$ws.Cells.Item(191, 6).Value = "Fruta"
$ws.Cells.Item(191, 7).Value = 100103
$ws.Cells.Item(191, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(191, 9).Value = 100103003
$ws.Cells.Item(191, 10).Value = "Damasco"
$ws.Cells.Item(191, 11).Value = "Helena"
$ws.Cells.Item(191, 12).Value = "Especial"
$ws.Cells.Item(191, 13).Value = 330
$ws.Cells.Item(191, 14).Value = 16000
$ws.Cells.Item(191, 15).Value = 16000
$ws.Cells.Item(191, 16).Value = 16000
$ws.Cells.Item(191, 17).Value = "$/caja 16 kilos"
$ws.Cells.Item(191, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(191, 19).Value = 1000
$ws.Cells.Item(191, 20).Value = 16

# Row 192
$ws.Cells.Item(192, 1).Value = 9
$ws.Cells.Item(192, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(192, 3).Value = "Metropolitana"
$ws.Cells.Item(192, 4).Value = 44558
$ws.Cells.Item(192, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(192, 5).Value = 13
$ws.Cells.Item(192, 6).Value = "Fruta"
$ws.Cells.Item(192, 7).Value = 100103
$ws.Cells.Item(192, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(192, 9).Value = 100103003
$ws.Cells.Item(192, 10).Value = "Damasco"
$ws.Cells.Item(192, 11).Value = "Patterson"
$ws.Cells.Item(192, 12).Value = "Primera"
$ws.Cells.Item(192, 13).Value = 410
$ws.Cells.Item(192, 14).Value = 14000
$ws.Cells.Item(192, 15).Value = 14000
$ws.Cells.Item(192, 16).Value = 14000
$ws.Cells.Item(192, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(192, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(192, 19).Value = 778
$ws.Cells.Item(192, 20).Value = 18

# Row 193
$ws.Cells.Item(193, 1).Value = 9
$ws.Cells.Item(193, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(193, 3).Value = "Metropolitana"
$ws.Cells.Item(193, 4).Value = 44558
$ws.Cells.Item(193, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(193, 5).Value = 13
$ws.Cells.Item(193, 6).Value = "Fruta"
$ws.Cells.Item(193, 7).Value = 100103
$ws.Cells.Item(193, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(193, 9).Value = 100103003
$ws.Cells.Item(193, 10).Value = "Damasco"
$ws.Cells.Item(193, 11).Value = "Patterson"
$ws.Cells.Item(193, 12).Value = "Segunda"
$ws.Cells.Item(193, 13).Value = 380
$ws.Cells.Item(193, 14).Value = 12000
$ws.Cells.Item(193, 15).Value = 12000
$ws.Cells.Item(193, 16).Value = 12000
$ws.Cells.Item(193, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(193, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(193, 19).Value = 667
$ws.Cells.Item(193, 20).Value = 18

# Row 194
$ws.Cells.Item(194, 1).Value = 9
$ws.Cells.Item(194, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(194, 3).Value = "Metropolitana"
$ws.Cells.Item(194, 4).Value = 44558
$ws.Cells.Item(194, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(194, 5).Value = 13
$ws.Cells.Item(194, 6).Value = "Fruta"
$ws.Cells.Item(194, 7).Value = 100103
$ws.Cells.Item(194, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(194, 9).Value = 100103003
$ws.Cells.Item(194, 10).Value = "Damasco"
$ws.Cells.Item(194, 11).Value = "Helena"
$ws.Cells.Item(194, 12).Value = "Primera"
$ws.Cells.Item(194, 13).Value = 350
$ws.Cells.Item(194, 14).Value = 12800
$ws.Cells.Item(194, 15).Value = 12800
$ws.Cells.Item(194, 16).Value = 12800
$ws.Cells.Item(194, 17).Value = "$/caja 16 kilos"
$ws.Cells.Item(194, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(194, 19).Value = 800
$ws.Cells.Item(194, 20).Value = 16

# Row 195
$ws.Cells.Item(195, 1).Value = 9
$ws.Cells.Item(195, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(195, 3).Value = "Metropolitana"
$ws.Cells.Item(195, 4).Value = 44957
$ws.Cells.Item(195, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(195, 5).Value = 13
$ws.Cells.Item(195, 6).Value = "Fruta"
$ws.Cells.Item(195, 7).Value = 100103
$ws.Cells.Item(195, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(195, 9).Value = 100103003
$ws.Cells.Item(195, 10).Value = "Damasco"
$ws.Cells.Item(195, 11).Value = "Modesto"
$ws.Cells.Item(195, 12).Value = "Primera"
$ws.Cells.Item(195, 13).Value = 250
$ws.Cells.Item(195, 14).Value = 9000
$ws.Cells.Item(195, 15).Value = 9000
$ws.Cells.Item(195, 16).Value = 9000
$ws.Cells.Item(195, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(195, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(195, 19).Value = 900
$ws.Cells.Item(195, 20).Value = 10

# Row 196
$ws.Cells.Item(196, 1).Value = 9
$ws.Cells.Item(196, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(196, 3).Value = "Metropolitana"
$ws.Cells.Item(196, 4).Value = 44957
$ws.Cells.Item(196, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(196, 5).Value = 13
$ws.Cells.Item(196, 6).Value = "Fruta"
$ws.Cells.Item(196, 7).Value = 100103
$ws.Cells.Item(196, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(196, 9).Value = 100103003
$ws.Cells.Item(196, 10).Value = "Damasco"
$ws.Cells.Item(196, 11).Value = "Modesto"
$ws.Cells.Item(196, 12).Value = "Segunda"
$ws.Cells.Item(196, 13).Value = 280
$ws.Cells.Item(196, 14).Value = 7000
$ws.Cells.Item(196, 15).Value = 7000
$ws.Cells.Item(196, 16).Value = 7000
$ws.Cells.Item(196, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(196, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(196, 19).Value = 700
$ws.Cells.Item(196, 20).Value = 10

# Row 197
$ws.Cells.Item(197, 1).Value = 9
$ws.Cells.Item(197, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(197, 3).Value = "Metropolitana"
$ws.Cells.Item(197, 4).Value = 44547
$ws.Cells.Item(197, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(197, 5).Value = 13
$ws.Cells.Item(197, 6).Value = "Fruta"
$ws.Cells.Item(197, 7).Value = 100103
$ws.Cells.Item(197, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(197, 9).Value = 100103003
$ws.Cells.Item(197, 10).Value = "Damasco"
$ws.Cells.Item(197, 11).Value = "Dina"
$ws.Cells.Item(197, 12).Value = "Especial"
$ws.Cells.Item(197, 13).Value = 450
$ws.Cells.Item(197, 14).Value = 12000
$ws.Cells.Item(197, 15).Value = 12000
$ws.Cells.Item(197, 16).Value = 12000
$ws.Cells.Item(197, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(197, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(197, 19).Value = 1200
$ws.Cells.Item(197, 20).Value = 10

# Row 198
$ws.Cells.Item(198, 1).Value = 9
$ws.Cells.Item(198, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(198, 3).Value = "Metropolitana"
$ws.Cells.Item(198, 4).Value = 44547
$ws.Cells.Item(198, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(198, 5).Value = 13
$ws.Cells.Item(198, 6).Value = "Fruta"
$ws.Cells.Item(198, 7).Value = 100103
$ws.Cells.Item(198, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(198, 9).Value = 100103003
$ws.Cells.Item(198, 10).Value = "Damasco"
$ws.Cells.Item(198, 11).Value = "Dina"
$ws.Cells.Item(198, 12).Value = "Primera"
$ws.Cells.Item(198, 13).Value = 410
$ws.Cells.Item(198, 14).Value = 10000
$ws.Cells.Item(198, 15).Value = 10000
$ws.Cells.Item(198, 16).Value = 10000
$ws.Cells.Item(198, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(198, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(198, 19).Value = 1000
$ws.Cells.Item(198, 20).Value = 10

# Row 199
$ws.Cells.Item(199, 1).Value = 9
$ws.Cells.Item(199, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(199, 3).Value = "Metropolitana"
$ws.Cells.Item(199, 4).Value = 44547
$ws.Cells.Item(199, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(199, 5).Value = 13
$ws.Cells.Item(199, 6).Value = "Fruta"
$ws.Cells.Item(199, 7).Value = 100103
$ws.Cells.Item(199, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(199, 9).Value = 100103003
$ws.Cells.Item(199, 10).Value = "Damasco"
$ws.Cells.Item(199, 11).Value = "Dina"
$ws.Cells.Item(199, 12).Value = "Segunda"
$ws.Cells.Item(199, 13).Value = 350
$ws.Cells.Item(199, 14).Value = 8000
$ws.Cells.Item(199, 15).Value = 8000
$ws.Cells.Item(199, 16).Value = 8000
$ws.Cells.Item(199, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(199, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(199, 19).Value = 800
$ws.Cells.Item(199, 20).Value = 10

# Row 200
$ws.Cells.Item(200, 1).Value = 9
$ws.Cells.Item(200, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(200, 3).Value = "Metropolitana"
$ws.Cells.Item(200, 4).Value = 44897
$ws.Cells.Item(200, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(200, 5).Value = 13
$ws.Cells.Item(200, 6).Value = "Fruta"
$ws.Cells.Item(200, 7).Value = 100103
$ws.Cells.Item(200, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(200, 9).Value = 100103003
$ws.Cells.Item(200, 10).Value = "Damasco"
$ws.Cells.Item(200, 11).Value = "Castle Brite"
$ws.Cells.Item(200, 12).Value = "Especial"
$ws.Cells.Item(200, 13).Value = 280
$ws.Cells.Item(200, 14).Value = 21600
$ws.Cells.Item(200, 15).Value = 21600
$ws.Cells.Item(200, 16).Value = 21600
$ws.Cells.Item(200, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(200, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(200, 19).Value = 1200
$ws.Cells.Item(200, 20).Value = 18

# Row 201
$ws.Cells.Item(201, 1).Value = 9
$ws.Cells.Item(201, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(201, 3).Value = "Metropolitana"
$ws.Cells.Item(201, 4).Value = 44897
$ws.Cells.Item(201, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(201, 5).Value = 13
$ws.Cells.Item(201, 6).Value = "Fruta"
$ws.Cells.Item(201, 7).Value = 100103
$ws.Cells.Item(201, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(201, 9).Value = 100103003
$ws.Cells.Item(201, 10).Value = "Damasco"
$ws.Cells.Item(201, 11).Value = "Castle Brite"
$ws.Cells.Item(201, 12).Value = "Primera"
$ws.Cells.Item(201, 13).Value = 250
$ws.Cells.Item(201, 14).Value = 18000
$ws.Cells.Item(201, 15).Value = 18000
$ws.Cells.Item(201, 16).Value = 18000
$ws.Cells.Item(201, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(201, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(201, 19).Value = 1000
$ws.Cells.Item(201, 20).Value = 18

# Row 202
$ws.Cells.Item(202, 1).Value = 9
$ws.Cells.Item(202, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(202, 3).Value = "Metropolitana"
$ws.Cells.Item(202, 4).Value = 44897
$ws.Cells.Item(202, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(202, 5).Value = 13
$ws.Cells.Item(202, 6).Value = "Fruta"
$ws.Cells.Item(202, 7).Value = 100103
$ws.Cells.Item(202, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(202, 9).Value = 100103003
$ws.Cells.Item(202, 10).Value = "Damasco"
$ws.Cells.Item(202, 11).Value = "Castle Brite"
$ws.Cells.Item(202, 12).Value = "Segunda"
$ws.Cells.Item(202, 13).Value = 210
$ws.Cells.Item(202, 14).Value = 14400
$ws.Cells.Item(202, 15).Value = 14400
$ws.Cells.Item(202, 16).Value = 14400
$ws.Cells.Item(202, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(202, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(202, 19).Value = 800
$ws.Cells.Item(202, 20).Value = 18

# Extend the declared sheet dimension to cover the 3 newly appended rows
$ws.Range("A1:T202") | Out-Null